$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.281.86'
$ws.Range("E2").Value = '  +0.76%  '
$ws.Range("D3").Value = '2.173.33'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''253.18'
$ws.Range("E5").Value = '  +5.91%  '
$ws.Range("E6").Value = '  -0.54%  '
$ws.Range("D7").Value = '''73.92'
$ws.Range("E7").Value = '  +1.15%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '''0.580'
$ws.Range("E9").Value = '  -1.19%  '
$ws.Range("D10").Value = '''40.59'
$ws.Range("E10").Value = '  +0.65%  '
$ws.Range("D11").Value = '''0.0910'
$ws.Range("E11").Value = '  -0.54%  '
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").Value = '''6.74'
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("D14").Value = '2.500.28'
$ws.Range("E14").Value = '  -0.48%  '
$ws.Range("D15").Value = '''14.12'
$ws.Range("E15").Value = '  -2.43%  '
$ws.Range("D16").Value = '2.168.05'
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("E17").Value = '  -2.84%  '
$ws.Range("D18").Value = '42.187.00'
$ws.Range("E18").Value = '  +0.86%  '
$ws.Range("D19").Value = '''0.0000102'
$ws.Range("E19").Value = '  -2.09%  '
$ws.Range("D20").Value = '''70.44'
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("D22").Value = '''226.28'
$ws.Range("E22").Value = '  -0.45%  '
$ws.Range("E23").Value = '  +3.19%  '
$ws.Range("D24").Value = '''9.43'
$ws.Range("E24").Value = '  -6.52%  '
$ws.Range("E25").Value = '  -0.12%  '
$ws.Range("D26").Value = '''10.40'
$ws.Range("E26").Value = '  -3.63%  '
$ws.Range("D27").Value = '''3.35'
$ws.Range("E27").Value = '  +2.18%  '
$ws.Range("E28").Value = '  -1.24%  '
$ws.Range("E29").Value = '  -1.87%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").Value = '''36.61'
$ws.Range("E30").Value = '  +10.69%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").Value = '''169.88'
$ws.Range("E31").Value = '  -1.21%  '
$ws.Range("D32").Value = '''19.95'
$ws.Range("E32").Value = '  -0.24%  '
$ws.Range("E33").Value = '  +3.41%  '
$ws.Range("E34").Value = '  -4.29%  '
$ws.Range("E35").Value = '  -0.64%  '
$ws.Range("E36").Value = '  +1.19%  '
$ws.Range("E37").Value = '  -3.54%  '
$ws.Range("D38").Value = '''0.0333'
$ws.Range("E38").Value = '  +6.47%  '
$ws.Range("E39").Value = '  -2.53%  '
$ws.Range("D40").Value = '''11.70'
$ws.Range("E40").Value = '  -4.83%  '
$ws.Range("E41").Value = '  +1.67%  '
$ws.Range("D42").Value = '''59.21'
$ws.Range("E42").Value = '  -0.80%  '
$ws.Range("D43").Value = '''5.13'
$ws.Range("E43").Value = '  -5.13%  '
$ws.Range("D44").Value = '''102.91'
$ws.Range("E44").Value = '  +4.72%  '
$ws.Range("D45").Value = '''0.469'
$ws.Range("E45").Value = '  +11.38%  '
$ws.Range("D46").Value = '''0.0970'
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("D47").Value = '''8.21'
$ws.Range("E47").Value = '  -3.37%  '
$ws.Range("E48").Value = '  +8.24%  '
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("E50").Value = '  -0.21%  '
$ws.Range("E51").Value = '  +0.34%  '
